# Applies the "Add files via upload" change to the Estatisticas sheet:
#  - Row 30 (C:F): convert existing text values "0.5","0","100","10000" into real numbers
#  - Rows 31-40: append 10 new test-result rows
#  - Row 40 (C:F): keep values as text (matching the row as committed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Entradas string reused by every appended row in this block
$entradas = "HTHG,HTAG,HTR,HS,AS,HST,AST,HC,AC,HY,AY,HR,AR"

# --- Row 30: re-type C30:F30 as numbers (values stay the same) ---
$ws.Cells.Item(30,3).Value = 0.5
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 100
$ws.Cells.Item(30,6).Value = 10000

# --- Rows 31-39: new rows, all numeric in columns C:F ---
$rows = @(
    @(31, "Teste Premier League 0.5 0,2", "E0", 0.5, 0.2, 100, 10000, "7 de 10"),
    @(32, "Teste", "D1", 0.5, 0.2, 100, 1000, "5 de 10"),
    @(33, "Teste", "D1", 0.5, 0.2, 100, 1000, "5 de 10"),
    @(34, "Teste", "D1", 0.5, 0.2, 100, 1000, "5 de 10"),
    @(35, "Teste", "D1", 0.5, 0.2, 100, 1000, "5 de 10"),
    @(36, "Teste", "D1", 0.5, 0.2, 100, 1000, "5 de 10"),
    @(37, "Teste", "D1", 0.5, 0.2, 100, 1000, "5 de 10"),
    @(38, "Teste", "D1", 0.5, 0.2, 100, 1000, "5 de 10"),
    @(39, "Teste", "E0", 0.5, 0.2, 100, 1000, "5 de 10")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r,1).Value = $row[1]
    $ws.Cells.Item($r,2).Value = $row[2]
    $ws.Cells.Item($r,3).Value = $row[3]
    $ws.Cells.Item($r,4).Value = $row[4]
    $ws.Cells.Item($r,5).Value = $row[5]
    $ws.Cells.Item($r,6).Value = $row[6]
    $ws.Cells.Item($r,7).Value = $row[7]
    $ws.Cells.Item($r,8).Value = $entradas
}

# --- Row 40: new row, C40:F40 stored as text (not numbers) ---
$ws.Cells.Item(40,1).Value = "Teste"
$ws.Cells.Item(40,2).Value = "E0"

$ws.Range("C40:F40").NumberFormat = "@"
$ws.Cells.Item(40,3).Value = "0.5"
$ws.Cells.Item(40,4).Value = "0.2"
$ws.Cells.Item(40,5).Value = "100"
$ws.Cells.Item(40,6).Value = "1000"

$ws.Cells.Item(40,7).Value = "6 de 10"
$ws.Cells.Item(40,8).Value = $entradas
